$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44305
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("Q2").Value = "$/caja 15 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1200

# Row 3
$ws.Range("D3").Value = 44305
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/caja 15 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44309
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = "$/caja 15 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44309
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/caja 15 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44285
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "$/caja 15 kilos empedrada"
$ws.Range("R6").Value = "Provincia del Elquí"
$ws.Range("S6").Value = 1200

# Row 7
$ws.Range("D7").Value = 44285
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 90
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = "$/caja 15 kilos empedrada"
$ws.Range("R7").Value = "Provincia del Elquí"
$ws.Range("S7").Value = 1000

# Row 8
$ws.Range("D8").Value = 44285
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 75
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = "$/caja 15 kilos empedrada"
$ws.Range("R8").Value = "Provincia del Elquí"
$ws.Range("S8").Value = 800
